$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Backend Developer"
$ws.Cells.Item(2, 2).Value = "Ensoft"
$ws.Cells.Item(2, 3).Value = "Link is not available"
$ws.Cells.Item(3, 1).Value = "Senior Software Developer"
$ws.Cells.Item(3, 2).Value = "Oracle"
$ws.Cells.Item(3, 3).Value = "https://id.indeed.com//cmp/Oracle"
$ws.Cells.Item(4, 1).Value = "Software Quality Assurance Intern"
$ws.Cells.Item(4, 2).Value = "Shopee"
$ws.Cells.Item(4, 3).Value = "https://id.indeed.com//cmp/Shopee"
$ws.Cells.Item(5, 1).Value = "Data Scientist"
$ws.Cells.Item(5, 2).Value = "Lancar"
$ws.Cells.Item(5, 3).Value = "Link is not available"
$ws.Cells.Item(6, 1).Value = "Lead Instructors - Le Wagon Data Science Bootcamp"
$ws.Cells.Item(6, 2).Value = "Le Wagon Bali"
$ws.Cells.Item(6, 3).Value = "Link is not available"
$ws.Cells.Item(7, 1).Value = "ERP Consultant / ERP Implementator"
$ws.Cells.Item(7, 2).Value = "HashMicro"
$ws.Cells.Item(7, 3).Value = "Link is not available"
$ws.Cells.Item(8, 1).Value = "DevOps Engineer (Work from Home; Full-time)"
$ws.Cells.Item(8, 2).Value = "Dynamic Technology Lab Pte Ltd"
$ws.Cells.Item(8, 3).Value = "Link is not available"
$ws.Cells.Item(9, 1).Value = "HRIS Developer"
$ws.Cells.Item(9, 2).Value = "Binabusana Internusa"
$ws.Cells.Item(9, 3).Value = "Link is not available"
$ws.Cells.Item(10, 1).Value = "Head of Data"
$ws.Cells.Item(10, 2).Value = "PT Sinar Mas Digital Ventures"
$ws.Cells.Item(10, 3).Value = "Link is not available"
$ws.Cells.Item(11, 1).Value = "Python Programmer"
$ws.Cells.Item(11, 2).Value = "1rstWAP"
$ws.Cells.Item(11, 3).Value = "Link is not available"
$ws.Cells.Item(12, 1).Value = "Publisher Support Specialist"
$ws.Cells.Item(12, 2).Value = "Coda Payments"
$ws.Cells.Item(12, 3).Value = "Link is not available"
$ws.Cells.Item(13, 1).Value = "Business Intelligence Developer E-Commerce"
$ws.Cells.Item(13, 2).Value = "Kompas Gramedia"
$ws.Cells.Item(13, 3).Value = "Link is not available"
$ws.Cells.Item(14, 1).Value = "ShopeePay QA Engineer [Entry Level] - NEW"
$ws.Cells.Item(14, 2).Value = "Shopee"
$ws.Cells.Item(14, 3).Value = "https://id.indeed.com//cmp/Shopee"
$ws.Cells.Item(15, 1).Value = "ShopeePay Backend Engineer [Experienced]"
$ws.Cells.Item(15, 2).Value = "Shopee"
$ws.Cells.Item(15, 3).Value = "https://id.indeed.com//cmp/Shopee"
$ws.Cells.Item(16, 1).Value = "Web Developer"
$ws.Cells.Item(16, 2).Value = "1rstWAP"
$ws.Cells.Item(16, 3).Value = "Link is not available"
$ws.Cells.Item(17, 1).Value = "ERP Developer"
$ws.Cells.Item(17, 2).Value = "PT Monotaro Indonesia"
$ws.Cells.Item(17, 3).Value = "Link is not available"
$ws.Cells.Item(18, 1).Value = "Senior Backend Developer"
$ws.Cells.Item(18, 2).Value = "AiChat Pte Ltd"
$ws.Cells.Item(18, 3).Value = "Link is not available"
$ws.Cells.Item(19, 1).Value = "Back End Developer"
$ws.Cells.Item(19, 2).Value = "Renos.id"
$ws.Cells.Item(19, 3).Value = "Link is not available"
$ws.Cells.Item(20, 1).Value = "PHP Developer"
$ws.Cells.Item(20, 2).Value = "PT Media Mitrakarya Indonesia"
$ws.Cells.Item(20, 3).Value = "Link is not available"
$ws.Cells.Item(21, 1).Value = "System Administrator"
$ws.Cells.Item(21, 2).Value = "Jawasoft"
$ws.Cells.Item(21, 3).Value = "Link is not available"
$ws.Cells.Item(22, 1).Value = "Full Stack/Backend Developer"
$ws.Cells.Item(22, 2).Value = "TPG Telecom Pte Ltd"
$ws.Cells.Item(22, 3).Value = "Link is not available"
$ws.Cells.Item(23, 1).Value = "Fullstack Developer"
$ws.Cells.Item(23, 2).Value = "PT Hermes Solusi Integrasi"
$ws.Cells.Item(23, 3).Value = "Link is not available"
$ws.Cells.Item(24, 1).Value = "Full Stack Developer"
$ws.Cells.Item(24, 2).Value = "PopBox"
$ws.Cells.Item(24, 3).Value = "Link is not available"
$ws.Cells.Item(25, 1).Value = "IT Production Support (Remote)"
$ws.Cells.Item(25, 2).Value = "mClinica"
$ws.Cells.Item(25, 3).Value = "Link is not available"
$ws.Cells.Item(26, 1).Value = "Full Stack Developer"
$ws.Cells.Item(26, 2).Value = "Ensoft"
$ws.Cells.Item(26, 3).Value = "Link is not available"
$ws.Cells.Item(27, 1).Value = "Web Developer"
$ws.Cells.Item(27, 2).Value = "Great Giant Foods"
$ws.Cells.Item(27, 3).Value = "Link is not available"
$ws.Cells.Item(28, 1).Value = "Software Quality Assurance - Manual Testing"
$ws.Cells.Item(28, 2).Value = "Cermati.com"
$ws.Cells.Item(28, 3).Value = "Link is not available"
$ws.Cells.Item(29, 1).Value = "ODOO developer/Junior&Senior Developer/Software engineer."
$ws.Cells.Item(29, 2).Value = "PT. Virgo Stellar"
$ws.Cells.Item(29, 3).Value = "Link is not available"
$ws.Cells.Item(30, 1).Value = "Application Developer"
$ws.Cells.Item(30, 2).Value = "Bank Mega"
$ws.Cells.Item(30, 3).Value = "https://id.indeed.com//cmp/Bank-Mega"
$ws.Cells.Item(31, 1).Value = "Golang Developer"
$ws.Cells.Item(31, 2).Value = "PT. Indocyber Global Technology"
$ws.Cells.Item(31, 3).Value = "Link is not available"
$ws.Cells.Item(32, 1).Value = "Developer / Programmer"
$ws.Cells.Item(32, 2).Value = "StrategArt"
$ws.Cells.Item(32, 3).Value = "Link is not available"
$ws.Cells.Item(33, 1).Value = "Web Developer"
$ws.Cells.Item(33, 2).Value = "Binabusana Internusa"
$ws.Cells.Item(33, 3).Value = "Link is not available"
$ws.Cells.Item(34, 1).Value = "Lead Software Quality Assurance"
$ws.Cells.Item(34, 2).Value = "Cermati.com"
$ws.Cells.Item(34, 3).Value = "Link is not available"
$ws.Cells.Item(35, 1).Value = "Back End Developer"
$ws.Cells.Item(35, 2).Value = "Akseleran"
$ws.Cells.Item(35, 3).Value = "Link is not available"
$ws.Cells.Item(36, 1).Value = "Developer"
$ws.Cells.Item(36, 2).Value = "Kinarya Alihdaya Mandiri PT"
$ws.Cells.Item(36, 3).Value = "Link is not available"
$ws.Cells.Item(37, 1).Value = "Senior Developer - Network (contract based)"
$ws.Cells.Item(37, 2).Value = "Standard Chartered"
$ws.Cells.Item(37, 3).Value = "Link is not available"
$ws.Cells.Item(38, 1).Value = "Technical Operations Engineer"
$ws.Cells.Item(38, 2).Value = "byOrange"
$ws.Cells.Item(38, 3).Value = "Link is not available"
$ws.Cells.Item(39, 1).Value = "Golang Developer (Back End)"
$ws.Cells.Item(39, 2).Value = "PT Lunaria Annua Teknologi (KoinWorks)"
$ws.Cells.Item(39, 3).Value = "Link is not available"
$ws.Cells.Item(40, 1).Value = "Datawarehouse Data Analyst (working in Kuala Lumpur, Malaysi..."
$ws.Cells.Item(40, 2).Value = "Mission Consultancy Services Malaysia SDN BHD"
$ws.Cells.Item(40, 3).Value = "Link is not available"
$ws.Cells.Item(41, 1).Value = "Front End Development – Consultant"
$ws.Cells.Item(41, 2).Value = "Accenture"
$ws.Cells.Item(41, 3).Value = "https://id.indeed.com//cmp/Accenture"
$ws.Cells.Item(42, 1).Value = "QA Engineer"
$ws.Cells.Item(42, 2).Value = "StyleTheory"
$ws.Cells.Item(42, 3).Value = "https://id.indeed.com//cmp/Styletheory"
$ws.Cells.Item(43, 1).Value = "IBM Service Associate Program - Application Developer"
$ws.Cells.Item(43, 2).Value = "IBM"
$ws.Cells.Item(43, 3).Value = "https://id.indeed.com//cmp/IBM"
$ws.Cells.Item(44, 1).Value = "Backend Engineer"
$ws.Cells.Item(44, 2).Value = "Cicil"
$ws.Cells.Item(44, 3).Value = "Link is not available"
$ws.Cells.Item(45, 1).Value = "IT WEB DEVELOPER"
$ws.Cells.Item(45, 2).Value = "Ismaya Group"
$ws.Cells.Item(45, 3).Value = "https://id.indeed.com//cmp/Ismaya-Group-1"
$ws.Cells.Item(46, 1).Value = "Associate Backend Engineer"
$ws.Cells.Item(46, 2).Value = "KeDA Tech"
$ws.Cells.Item(46, 3).Value = "Link is not available"
$ws.Cells.Item(47, 1).Value = "Data Warehouse Engineer"
$ws.Cells.Item(47, 2).Value = "GO-JEK"
$ws.Cells.Item(47, 3).Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Cells.Item(48, 1).Value = "Data Warehouse Engineer - GoPay"
$ws.Cells.Item(48, 2).Value = "GO-JEK"
$ws.Cells.Item(48, 3).Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Cells.Item(49, 1).Value = "ShopeePay QA Engineer [Experienced]"
$ws.Cells.Item(49, 2).Value = "Shopee"
$ws.Cells.Item(49, 3).Value = "https://id.indeed.com//cmp/Shopee"
$ws.Cells.Item(50, 1).Value = "Java Developer"
$ws.Cells.Item(50, 2).Value = "NTT Ltd"
$ws.Cells.Item(50, 3).Value = "Link is not available"
$ws.Cells.Item(51, 1).Value = "Full Stack Developer"
$ws.Cells.Item(51, 2).Value = "Vicuna Corp"
$ws.Cells.Item(51, 3).Value = "Link is not available"
$ws.Cells.Item(52, 1).Value = "Senior Developer"
$ws.Cells.Item(52, 2).Value = "PT Chrombit Digtal Lab"
$ws.Cells.Item(52, 3).Value = "Link is not available"
$ws.Cells.Item(53, 1).Value = "System Engineering Development"
$ws.Cells.Item(53, 2).Value = "Alodokter"
$ws.Cells.Item(53, 3).Value = "Link is not available"
$ws.Cells.Item(54, 1).Value = "Financial Service SRE Engineer [Entry Level]"
$ws.Cells.Item(54, 2).Value = "Shopee"
$ws.Cells.Item(54, 3).Value = "https://id.indeed.com//cmp/Shopee"
$ws.Cells.Item(55, 1).Value = "Business Intelligence Developer"
$ws.Cells.Item(55, 2).Value = "Stockbit-Bibit"
$ws.Cells.Item(55, 3).Value = "Link is not available"
$ws.Cells.Item(56, 1).Value = "ERP Programmer (Odoo Framework)"
$ws.Cells.Item(56, 2).Value = "HashMicro"
$ws.Cells.Item(56, 3).Value = "Link is not available"
$ws.Cells.Item(57, 1).Value = "Python Developer"
$ws.Cells.Item(57, 2).Value = "QSI Recruitment"
$ws.Cells.Item(57, 3).Value = "Link is not available"
$ws.Cells.Item(58, 1).Value = "Software Engineer - Data Platform"
$ws.Cells.Item(58, 2).Value = "Cermati.com"
$ws.Cells.Item(58, 3).Value = "Link is not available"
$ws.Cells.Item(59, 1).Value = "Senior Backend Developer"
$ws.Cells.Item(59, 2).Value = "Schoters"
$ws.Cells.Item(59, 3).Value = "Link is not available"
$ws.Cells.Item(60, 1).Value = "Senior Machine Learning"
$ws.Cells.Item(60, 2).Value = "Alodokter"
$ws.Cells.Item(60, 3).Value = "Link is not available"
$ws.Cells.Item(61, 1).Value = "Integration Developer"
$ws.Cells.Item(61, 2).Value = "GO-JEK"
$ws.Cells.Item(61, 3).Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Cells.Item(62, 1).Value = "Customer Solutions Consultant, Infrastructure Modernization,..."
$ws.Cells.Item(62, 2).Value = "Google"
$ws.Cells.Item(62, 3).Value = "https://id.indeed.com//cmp/Google"
$ws.Cells.Item(63, 1).Value = "Software Engineer, xShop"
$ws.Cells.Item(63, 2).Value = "Coda Payments"
$ws.Cells.Item(63, 3).Value = "Link is not available"
$ws.Cells.Item(64, 1).Value = "QA Engineer - GoFinance"
$ws.Cells.Item(64, 2).Value = "GO-JEK"
$ws.Cells.Item(64, 3).Value = "https://id.indeed.com//cmp/Pt.-Go--jek-Indonesia-2"
$ws.Cells.Item(65, 1).Value = "ShopeePay Backend Engineer [Leader]"
$ws.Cells.Item(65, 2).Value = "Shopee"
$ws.Cells.Item(65, 3).Value = "https://id.indeed.com//cmp/Shopee"
$ws.Cells.Item(66, 1).Value = "API Developer"
$ws.Cells.Item(66, 2).Value = "PT Multi Bangun Abadi"
$ws.Cells.Item(66, 3).Value = "https://id.indeed.com//cmp/PT-Multi-Bangun-Abadi"
$ws.Cells.Item(67, 1).Value = "Software Engineer"
$ws.Cells.Item(67, 2).Value = "Alterra"
$ws.Cells.Item(67, 3).Value = "Link is not available"
$ws.Cells.Item(68, 1).Value = "Data Engineer"
$ws.Cells.Item(68, 2).Value = "Tokenomy"
$ws.Cells.Item(68, 3).Value = "Link is not available"
$ws.Cells.Item(69, 1).Value = "BACKEND DEVELOPER"
$ws.Cells.Item(69, 2).Value = "Dipstrategy"
$ws.Cells.Item(69, 3).Value = "Link is not available"
$ws.Cells.Item(70, 1).Value = "Engineering and Technology - Sea Labs - Back End Engineer, P..."
$ws.Cells.Item(70, 2).Value = "Shopee"
$ws.Cells.Item(70, 3).Value = "https://id.indeed.com//cmp/Shopee"
$ws.Cells.Item(71, 1).Value = "Senior Back End Developer"
$ws.Cells.Item(71, 2).Value = "Sonar Social Media Monitoring Platform"
$ws.Cells.Item(71, 3).Value = "Link is not available"
$ws.Cells.Item(72, 1).Value = "Test Engineer"
$ws.Cells.Item(72, 2).Value = "Quipper"
$ws.Cells.Item(72, 3).Value = "Link is not available"
$ws.Cells.Item(73, 1).Value = "Productivity Engineer"
$ws.Cells.Item(73, 2).Value = "Stockbit-Bibit"
$ws.Cells.Item(73, 3).Value = "Link is not available"
# Row 74 (Principal Engineer / Stockbit / Link is not available) is unchanged
$ws.Cells.Item(75, 1).Value = "Engineering and Technology - Sea Labs - System Quality Assur..."
$ws.Cells.Item(75, 2).Value = "Shopee"
$ws.Cells.Item(75, 3).Value = "https://id.indeed.com//cmp/Shopee"
$ws.Cells.Item(76, 1).Value = "QA Automation Engineer"
$ws.Cells.Item(76, 2).Value = "Nimbly"
$ws.Cells.Item(76, 3).Value = "Link is not available"
